# pq finished and fixed reading time bug
#
# 1) Add a new row (31) for a new questionnaire string pair:
#    pq_check_and_change_answers / "Please check all your answers carefully.
#    If you want to change an answer, tick the box next to it." — styled
#    with an explicit black font (matches the new font/cellXf added upstream).
# 2) Shorten the pq_final_message text in B21 (drop the trailing
#    Submit/Previous sentence, now covered by the new check-your-answers
#    message instead).
# 3) Leave the selection on B21 (the cell that was last edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "pq_check_and_change_answers"
$ws.Range("B31").Value = "Please check all your answers carefully. If you want to change an answer, tick the box next to it."
$ws.Range("A31:B31").Font.Color = 0

$ws.Range("B21").Value = "Thank you for filling out our questionnaire. Remember that all data will be anonymised and treated confidentially."

$ws.Range("B21").Select() | Out-Null
